$wb = $excel.ActiveWorkbook

$errorMsg = 'The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/262b53c48ef1d08f5b9c2ce038e0b87fab1686cb/e2e/63f99f58-911e-434e-8014-f3d346862426.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/41724c0408642e8bba0767ff68543a34acb42285/e2e/63f99f58-911e-434e-8014-f3d346862426.md.'
$currentMdUrl = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/262b53c48ef1d08f5b9c2ce038e0b87fab1686cb/e2e/63f99f58-911e-434e-8014-f3d346862426.md'
$mdDisplay = '63f99f58-911e-434e-8014-f3d346862426.md'

# ---------- zh-cn sheet ----------
$ws = $wb.Worksheets.Item("zh-cn")

# Widen the "Error Detail" column (P) to fit the long message.
$ws.Columns.Item(16).ColumnWidth = 39.17

# Row 7 ("63f99f58-911e-434e-8014-f3d346862426") now has a target/handback
# file recorded, but it is not the latest version, so the datetime and
# error-detail columns get filled in too.
$ws.Hyperlinks.Add($ws.Range("I7"), $currentMdUrl, "", "", $mdDisplay)
$ws.Range("J7").Value2 = $ws.Range("G7").Value2
$ws.Range("K7").Value2 = "2016-08-31 18:49:49"
$ws.Range("P7").Value2 = $errorMsg

# ---------- de-de sheet ----------
$ws = $wb.Worksheets.Item("de-de")

$ws.Columns.Item(16).ColumnWidth = 39.17

$ws.Hyperlinks.Add($ws.Range("I7"), $currentMdUrl, "", "", $mdDisplay)
$ws.Range("J7").Value2 = $ws.Range("G7").Value2
$ws.Range("K7").Value2 = "2016-08-31 18:49:56"
$ws.Range("P7").Value2 = $errorMsg
